# Update attendance ("想去人数") figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of worksheet name -> hashtable of cell address -> new value
$updates = @{
    "展览" = @{
        "F2"  = 6368
        "F4"  = 4
        "F5"  = 373
        "F9"  = 69
        "F10" = 69
        "F12" = 150
        "F13" = 362
        "F14" = 620
        "F15" = 3104
        "F18" = 1760
        "F19" = 20
    }
    "全部类型" = @{
        "F2"  = 6368
        "F4"  = 4
        "F5"  = 373
        "F10" = 69
        "F11" = 69
        "F13" = 150
        "F14" = 362
        "F15" = 620
        "F16" = 3104
        "F19" = 1760
        "F20" = 20
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range($addr).Value = $cellUpdates[$addr]
    }
}
